$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 524
$newTimestamp = "2023-01-16 12:59:29"

# Refresh the crawl timestamp in column O (timestamp) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Updated ratingAmount (column D) for the Max Havelaar Mango product (row 83): 1 -> 2
$ws.Cells.Item(83, 4).Value = 2
